$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "season record" header columns, using the same header formatting
# (bold font, thin border, centered) as the existing header cells.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in each player's team season record (Wins/Losses/Ties) for every
# data row.
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 67   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 95   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
